$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1622.0834
$ws.Range("I19").Value = 1211.25
$ws.Range("J19").Value = 1827.5
$ws.Range("K19").Value = 1211.25
$ws.Range("L19").Value = 1827.5
$ws.Range("M19").Value = -1036.25
$ws.Range("N19").Value = -2177.5
$ws.Range("H80").Value = 4412.8115
$ws.Range("I80").Value = 3415.6365
$ws.Range("K80").Value = 10246.9095
$ws.Range("M80").Value = -9248.9095
$ws.Range("H83").Value = 4412.8115
$ws.Range("I83").Value = 3415.6365
$ws.Range("K83").Value = 30740.7285
$ws.Range("M83").Value = -25748.7285
$ws.Range("H86").Value = 6327.9414
$ws.Range("I86").Value = 4714.8184
$ws.Range("J86").Value = 9285.333000000001
$ws.Range("K86").Value = 4714.8184
$ws.Range("L86").Value = 9285.333000000001
$ws.Range("M86").Value = -3591.8184
$ws.Range("N86").Value = -11531.333
$ws.Range("H88").Value = 2991.3215
$ws.Range("I88").Value = 2031.2222
$ws.Range("J88").Value = 3446.1052
$ws.Range("K88").Value = 2031.2222
$ws.Range("L88").Value = 3446.1052
$ws.Range("M88").Value = -1625.2222
$ws.Range("N88").Value = -4258.1052
$ws.Range("H89").Value = 6327.9414
$ws.Range("I89").Value = 4714.8184
$ws.Range("J89").Value = 9285.333000000001
$ws.Range("K89").Value = 23574.092
$ws.Range("L89").Value = 46426.665
$ws.Range("M89").Value = -17958.092
$ws.Range("N89").Value = -57658.665
$ws.Range("H91").Value = 2991.3215
$ws.Range("I91").Value = 2031.2222
$ws.Range("J91").Value = 3446.1052
$ws.Range("K91").Value = 2031.2222
$ws.Range("L91").Value = 3446.1052
$ws.Range("M91").Value = -627.2221999999999
$ws.Range("N91").Value = -6254.1052
$ws.Range("H137").Value = 3040.4827
$ws.Range("I137").Value = 2473.375
$ws.Range("J137").Value = 3738.4614
$ws.Range("K137").Value = 7420.125
$ws.Range("L137").Value = 11215.3842
$ws.Range("M137").Value = -4870.125
$ws.Range("N137").Value = -16315.3842

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1158.2106
$ws.Range("I61").Value = 1087.4667
$ws.Range("J61").Value = 1423.5
$ws.Range("K61").Value = 1087.4667
$ws.Range("L61").Value = 1423.5
$ws.Range("M61").Value = -875.4666999999999
$ws.Range("N61").Value = -1847.5
$ws.Range("H74").Value = 42767.668
$ws.Range("I74").Value = 72314.71000000001
$ws.Range("J74").Value = 1401.8
$ws.Range("K74").Value = 72314.71000000001
$ws.Range("L74").Value = 1401.8
$ws.Range("M74").Value = -71440.71000000001
$ws.Range("N74").Value = -3149.8
$ws.Range("H77").Value = 42767.668
$ws.Range("I77").Value = 72314.71000000001
$ws.Range("J77").Value = 1401.8
$ws.Range("K77").Value = 361573.55
$ws.Range("L77").Value = 7009
$ws.Range("M77").Value = -357205.55
$ws.Range("N77").Value = -15745
$ws.Range("H122").Value = 2010.7241
$ws.Range("I122").Value = 1421.0714
$ws.Range("J122").Value = 2561.0667
$ws.Range("K122").Value = 4263.2142
$ws.Range("L122").Value = 7683.2001
$ws.Range("M122").Value = -1813.2142
$ws.Range("N122").Value = -12583.2001
$ws.Range("H132").Value = 2539.4062
$ws.Range("I132").Value = 2377.7083
$ws.Range("J132").Value = 3024.5
$ws.Range("K132").Value = 7133.124899999999
$ws.Range("L132").Value = 9073.5
$ws.Range("M132").Value = -4603.124899999999
$ws.Range("N132").Value = -14133.5
$ws.Range("H136").Value = 1158.2106
$ws.Range("I136").Value = 1087.4667
$ws.Range("J136").Value = 1423.5
$ws.Range("K136").Value = 3262.4001
$ws.Range("L136").Value = 4270.5
$ws.Range("M136").Value = -712.4000999999998
$ws.Range("N136").Value = -9370.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1745.1666
$ws.Range("I86").Value = 1810.6
$ws.Range("J86").Value = 1663.375
$ws.Range("K86").Value = 1810.6
$ws.Range("L86").Value = 1663.375
$ws.Range("M86").Value = -687.5999999999999
$ws.Range("N86").Value = -3909.375
$ws.Range("H89").Value = 1745.1666
$ws.Range("I89").Value = 1810.6
$ws.Range("J89").Value = 1663.375
$ws.Range("K89").Value = 9053
$ws.Range("L89").Value = 8316.875
$ws.Range("M89").Value = -3437
$ws.Range("N89").Value = -19548.875
$ws.Range("H99").Value = 1676.9131
$ws.Range("I99").Value = 1272.6
$ws.Range("J99").Value = 2435
$ws.Range("K99").Value = 1272.6
$ws.Range("L99").Value = 2435
$ws.Range("M99").Value = 225.4000000000001
$ws.Range("N99").Value = -5431
$ws.Range("H134").Value = 2969.6667
$ws.Range("I134").Value = 2577.6553
$ws.Range("J134").Value = 5811.75
$ws.Range("K134").Value = 7732.965899999999
$ws.Range("L134").Value = 17435.25
$ws.Range("M134").Value = -5197.965899999999
$ws.Range("N134").Value = -22505.25

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 39474908
$ws.Range("I31").Value = 35715230
$ws.Range("J31").Value = 50002000
$ws.Range("K31").Value = 35715230
$ws.Range("L31").Value = 50002000
$ws.Range("M31").Value = -35714935
$ws.Range("N31").Value = -50002590
$ws.Range("H34").Value = 39474908
$ws.Range("I34").Value = 35715230
$ws.Range("J34").Value = 50002000
$ws.Range("K34").Value = 35715230
$ws.Range("L34").Value = 50002000
$ws.Range("M34").Value = -35715028
$ws.Range("N34").Value = -50002404
$ws.Range("H58").Value = 4467.8
$ws.Range("I58").Value = 4847.154
$ws.Range("K58").Value = 4847.154
$ws.Range("M58").Value = -4644.154
$ws.Range("H99").Value = 4250
$ws.Range("I99").Value = 4550
$ws.Range("J99").Value = 4100
$ws.Range("K99").Value = 4550
$ws.Range("L99").Value = 4100
$ws.Range("M99").Value = -3052
$ws.Range("N99").Value = -7096
$ws.Range("H126").Value = 4250
$ws.Range("I126").Value = 4550
$ws.Range("J126").Value = 4100
$ws.Range("K126").Value = 13650
$ws.Range("L126").Value = 12300
$ws.Range("M126").Value = -11180
$ws.Range("N126").Value = -17240
$ws.Range("H132").Value = 1964.122
$ws.Range("I132").Value = 1646.2
$ws.Range("K132").Value = 4938.6
$ws.Range("M132").Value = -2408.6
$ws.Range("H134").Value = 1682
$ws.Range("I134").Value = 1727.5333
$ws.Range("J134").Value = 999
$ws.Range("K134").Value = 5182.5999
$ws.Range("L134").Value = 2997
$ws.Range("M134").Value = -2647.5999
$ws.Range("N134").Value = -8067
$ws.Range("H136").Value = 4467.8
$ws.Range("I136").Value = 4847.154
$ws.Range("K136").Value = 14541.462
$ws.Range("M136").Value = -11991.462

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 567.1795
$ws.Range("I113").Value = 585.3913
$ws.Range("J113").Value = 541
$ws.Range("K113").Value = 1756.1739
$ws.Range("L113").Value = 1623
$ws.Range("M113").Value = 413.8261
$ws.Range("N113").Value = -5963
$ws.Range("H131").Value = 1005.425
$ws.Range("I131").Value = 821.6
$ws.Range("J131").Value = 1066.7
$ws.Range("K131").Value = 2464.8
$ws.Range("L131").Value = 3200.1
$ws.Range("M131").Value = 2575.2
$ws.Range("N131").Value = -13280.1

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2700
$ws.Range("I102").Value = 2566.6667
$ws.Range("J102").Value = 2900
$ws.Range("K102").Value = 2566.6667
$ws.Range("L102").Value = 2900
$ws.Range("M102").Value = -944.6667000000002
$ws.Range("N102").Value = -6144
$ws.Range("H132").Value = 2820.25
$ws.Range("I132").Value = 2230.625
$ws.Range("J132").Value = 3999.5
$ws.Range("K132").Value = 6691.875
$ws.Range("L132").Value = 11998.5
$ws.Range("M132").Value = -4161.875
$ws.Range("N132").Value = -17058.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3171.5715
$ws.Range("I132").Value = 2138
$ws.Range("J132").Value = 4549.6665
$ws.Range("K132").Value = 6414
$ws.Range("L132").Value = 13648.9995
$ws.Range("M132").Value = -3884
$ws.Range("N132").Value = -18708.9995
$ws.Range("H136").Value = 12822056
$ws.Range("I136").Value = 23810704
$ws.Range("J136").Value = 1967.5
$ws.Range("K136").Value = 71432112
$ws.Range("L136").Value = 5902.5
$ws.Range("M136").Value = -71429562
$ws.Range("N136").Value = -11002.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 6898.3335
$ws.Range("I29").Value = 5347.5
$ws.Range("K29").Value = 5347.5
$ws.Range("M29").Value = -5057.5
$ws.Range("H49").Value = 8000
$ws.Range("J49").Value = 8000
$ws.Range("L49").Value = 8000
$ws.Range("N49").Value = -8460
$ws.Range("H122").Value = 13222980
$ws.Range("I122").Value = 31251488
$ws.Range("J122").Value = 111337.73
$ws.Range("K122").Value = 93754464
$ws.Range("L122").Value = 334013.19
$ws.Range("M122").Value = -93752014
$ws.Range("N122").Value = -338913.19
$ws.Range("H132").Value = 4765202.5
$ws.Range("I132").Value = 6253358
$ws.Range("J132").Value = 3105.7
$ws.Range("K132").Value = 18760074
$ws.Range("L132").Value = 9317.099999999999
$ws.Range("M132").Value = -18757544
$ws.Range("N132").Value = -14377.1
$ws.Range("H136").Value = 34484490
$ws.Range("I136").Value = 66667950
$ws.Range("J136").Value = 2204.1428
$ws.Range("K136").Value = 200003850
$ws.Range("L136").Value = 6612.428400000001
$ws.Range("M136").Value = -200001300
$ws.Range("N136").Value = -11712.4284
